$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "description" column (D) for the rows that were missing it
# (rows 7-12) with the new popover/description text, matching the commit's
# "tweaked the film popover" text entry, then leave the range selected.
$rng = $ws.Range("D7:D12")
$rng.Value = "a good film"
$rng.Select()
